$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data to append (row 27)
$newDateSerial = 46022   # 31-Dec-2025 in Excel's 1900 date system
$newText = "animation for the background and finished crow arc"
$newHours = 6

# Copy the style of the previous date cell (A26) into A27 so the new date
# gets the same date number format, then set the value.
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A27").Value = $newDateSerial
$ws.Range("B27").Value = $newText
$ws.Range("C27").Value = $newHours

# Update the selection to mirror the authored workbook view state.
$ws.Range("E26").Select() | Out-Null

$wb.Save()
